# Updated cryptos list on Tue Jan 30 15:35:23 UTC 2024 with GitHub Actions
# Applies the refreshed price / 1h-volume figures (and a couple of rank
# swaps between rows that exchanged positions) to the crypto tracker sheet.
# Price cells are forced to text (NumberFormat "@") before assignment, then
# restored to the default "Normal" style so numeric-looking strings such as
# "105.86" are preserved verbatim instead of being coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.473.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.315.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.90%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '105.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.45%  '
$ws.Range("E10").Value = '  +3.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0814'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.670.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.352.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.803'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.444.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0924'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '241.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.32%  '
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.26%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.53%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.85%  '
$ws.Range("E36").Value = '  +6.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0737'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.53%  '
$ws.Range("E38").Value = '  +13.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.00%  '
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0291'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.964.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("E47").Value = '  +6.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.09%  '
